$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its values as text, matching the
# source data which stores prices as literal strings (not numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.243.82"
$ws.Range("E2").Value = "  +2.47%  "

$ws.Range("D3").Value = "1.876.98"
$ws.Range("E3").Value = "  +4.81%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").Value = "311.66"
$ws.Range("E5").Value = "  +2.29%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").Value = "0.5072"
$ws.Range("E7").Value = "  +2.87%  "

$ws.Range("D8").Value = "0.3930"
$ws.Range("E8").Value = "  +3.09%  "

$ws.Range("D9").Value = "0.09841"
$ws.Range("E9").Value = "  +7.46%  "

$ws.Range("D10").Value = "1.148"
$ws.Range("E10").Value = "  +5.66%  "

$ws.Range("D11").Value = "40.89"
$ws.Range("E11").Value = "  +1.09%  "

$ws.Range("D12").Value = "6.504"
$ws.Range("E12").Value = "  +4.63%  "

$ws.Range("D13").Value = "21.04"
$ws.Range("E13").Value = "  +3.68%  "

$ws.Range("D14").Value = "1.878.09"
$ws.Range("E14").Value = "  +4.87%  "

$ws.Range("D15").Value = "7.456"
$ws.Range("E15").Value = "  +4.85%  "

$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").Value = "0.00001133"
$ws.Range("E17").Value = "  +3.38%  "

$ws.Range("D18").Value = "93.14"
$ws.Range("E18").Value = "  +1.63%  "

$ws.Range("D19").Value = "0.06595"
$ws.Range("E19").Value = "  +0.66%  "

$ws.Range("D20").Value = "17.65"
$ws.Range("E20").Value = "  +4.10%  "

$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").Value = "6.201"
$ws.Range("E22").Value = "  +5.03%  "

$ws.Range("D23").Value = "28.301.13"
$ws.Range("E23").Value = "  +2.51%  "

$ws.Range("D24").Value = "11.34"
$ws.Range("E24").Value = "  +3.72%  "

$ws.Range("D25").Value = "2.296"
$ws.Range("E25").Value = "  +4.01%  "

$ws.Range("D26").Value = "2.587"
$ws.Range("E26").Value = "  +9.60%  "

$ws.Range("D27").Value = "2.094.99"
$ws.Range("E27").Value = "  +4.83%  "

$ws.Range("E28").Value = "  +5.25%  "

$ws.Range("D29").Value = "158.98"
$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("D30").Value = "127.74"
$ws.Range("E30").Value = "  +1.26%  "

$ws.Range("D31").Value = "0.1069"
$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").Value = "1.077"
$ws.Range("E32").Value = "  +3.38%  "

$ws.Range("D33").Value = "5.654"
$ws.Range("E33").Value = "  +3.27%  "

$ws.Range("D34").Value = "3.622"
$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").Value = "9.547"
$ws.Range("E35").Value = "  +9.16%  "

$ws.Range("D36").Value = "0.06740"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").Value = "0.02388"
$ws.Range("E37").Value = "  +4.86%  "

$ws.Range("D38").Value = "0.2206"
$ws.Range("E38").Value = "  +4.50%  "

$ws.Range("D39").Value = "0.6392"
$ws.Range("E39").Value = "  +5.24%  "

$ws.Range("D40").Value = "11.53"
$ws.Range("E40").Value = "  +2.63%  "

$ws.Range("D41").Value = "5.008"
$ws.Range("E41").Value = "  +2.81%  "

$ws.Range("E42").Value = "  +4.53%  "

$ws.Range("D43").Value = "0.9997"
$ws.Range("E43").Value = "  -0.65%  "

$ws.Range("D44").Value = "13.53"
$ws.Range("E44").Value = "  +4.18%  "

$ws.Range("D45").Value = "0.6014"
$ws.Range("E45").Value = "  +4.05%  "

$ws.Range("D46").Value = "3.662"
$ws.Range("E46").Value = "  +0.43%  "

$ws.Range("D47").Value = "1.270"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("E48").Value = "  +5.01%  "

$ws.Range("D49").Value = "124.40"
$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("D50").Value = "1.199"
$ws.Range("E50").Value = "  +3.53%  "

$ws.Range("D51").Value = "0.06860"
$ws.Range("E51").Value = "  +2.39%  "
